# The "sp-2" and "sp-3" tabs swap names (the underlying sheet data/IDs stay
# put - only the display names trade places), and the previously-active
# tab ("sp-1") is replaced by the sheet that ends up named "sp-3" (the one
# that used to be called "sp-2").
$wb = $excel.ActiveWorkbook

$sheetB = $wb.Worksheets.Item("sp-2")   # will become "sp-3"
$sheetC = $wb.Worksheets.Item("sp-3")   # will become "sp-2"

# Use a scratch name so the two final names never collide mid-swap.
$sheetB.Name = "sp-2-swap-tmp"
$sheetC.Name = "sp-2"
$sheetB.Name = "sp-3"

# Make the (renamed) "sp-3" sheet the active / selected tab, moving
# tabSelected + activeTab off of "sp-1" and onto it.
$wb.Worksheets.Item("sp-3").Activate()
